$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.437.06'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '2.240.75'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.833'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").Value = '2.213.54'
$ws.Range("E15").Value = '  -2.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '44.082.49'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '0.0₃0954'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.72%  '
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").Value = '  +4.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0792'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("E33").Value = '  +2.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.110'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.07%  '
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0297'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").Value = '1.783.70'
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("E44").Value = '  +2.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '78.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '70.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.14%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.51%  '
